$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.896.21'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '2.032.91'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.96'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.638'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '62.68'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.396'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +7.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '57.93'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.61%  '
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.896'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.46'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +18.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.44'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').Value = '2.329.59'
$ws.Range('E16').Value = '  +0.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.59'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.90%  '
$ws.Range('D18').Value = '2.030.44'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').Value = '36.810.03'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.48'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('E21').Value = '  +2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.39'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.69'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.53'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.24%  '
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.90'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.141'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +20.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '160.65'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.36'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.84%  '
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.19'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.05'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0625'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.17%  '
$ws.Range('E35').Value = '  +0.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.59'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +11.53%  '
$ws.Range('E37').Value = '  -3.79%  '
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.19'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +23.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.102'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.27'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +3.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.93'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0216'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.97'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '94.14'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.72'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('D49').Value = '1.365.40'
$ws.Range('E49').Value = '  -4.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.92'
$ws.Range('D50').ClearFormats()
$ws.Range('D51').Value = '2.224.75'
$ws.Range('E51').Value = '  +1.33%  '
